# Apply updated TPM-derived NATMI metrics to the LR-pair worksheet.
# Values below are the recomputed outputs (ligand/receptor detection,
# expression, specificity and edge-weight columns) after refreshing the
# underlying TPM data; cell/row layout is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.699658666666667
$ws.Range("N2").Value = 20.098976
$ws.Range("O2").Value = 0.1402150605386345
$ws.Range("P2").Value = 0.1402150605386345
$ws.Range("Q2").Value = 3.057751014101334
$ws.Range("R2").Value = 27.519759126912
$ws.Range("S2").Value = 0.002685063398746761
$ws.Range("T2").Value = 0.002685063398746761

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.6453289538613627
$ws.Range("P3").Value = 0.6453289538613627
$ws.Range("Q3").Value = 14.07306216264
$ws.Range("R3").Value = 126.65755946376
$ws.Range("S3").Value = 0.01235779628456706
$ws.Range("T3").Value = 0.01235779628456706

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.2144559856000028
$ws.Range("P4").Value = 0.2144559856000028
$ws.Range("Q4").Value = 4.676765854748
$ws.Range("R4").Value = 42.090892692732
$ws.Range("S4").Value = 0.004106747986733336
$ws.Range("T4").Value = 0.004106747986733336

# Row 5
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.699658666666667
$ws.Range("N5").Value = 20.098976
$ws.Range("O5").Value = 0.1402150605386345
$ws.Range("P5").Value = 0.1402150605386345
$ws.Range("Q5").Value = 132.2927485996836
$ws.Range("R5").Value = 1190.634737397152
$ws.Range("S5").Value = 0.1161685224030623
$ws.Range("T5").Value = 0.1161685224030623

# Row 6
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.6453289538613627
$ws.Range("P6").Value = 0.6453289538613627
$ws.Range("S6").Value = 0.5346566249445955
$ws.Range("T6").Value = 0.5346566249445955

# Row 7
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.2144559856000028
$ws.Range("P7").Value = 0.2144559856000028
$ws.Range("Q7").Value = 202.3389761391663
$ws.Range("S7").Value = 0.1776773113525865
$ws.Range("T7").Value = 0.1776773113525865

# Row 8
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.699658666666667
$ws.Range("N8").Value = 20.098976
$ws.Range("O8").Value = 0.1402150605386345
$ws.Range("P8").Value = 0.1402150605386345
$ws.Range("Q8").Value = 24.32645391900801
$ws.Range("R8").Value = 218.938085271072
$ws.Range("S8").Value = 0.02136147473682547
$ws.Range("T8").Value = 0.02136147473682547

# Row 9
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.6453289538613627
$ws.Range("P9").Value = 0.6453289538613627
$ws.Range("S9").Value = 0.09831453263220019
$ws.Range("T9").Value = 0.09831453263220019

# Row 10
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.2144559856000028
$ws.Range("P10").Value = 0.2144559856000028
$ws.Range("Q10").Value = 37.206799550013
$ws.Range("R10").Value = 334.861195950117
$ws.Range("S10").Value = 0.03267192626068299
$ws.Range("T10").Value = 0.03267192626068298
